$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.224.58'
$ws.Range("E2").Value = '  +1.23%  '
$ws.Range("D3").Value = '1.875.61'
$ws.Range("E3").Value = '  +3.21%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '311.62'
$ws.Range("E5").Value = '  +0.30%  '
$ws.Range("E6").Value = '  +0.05%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5010'
$ws.Range("E7").Value = '  +0.60%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3893'
$ws.Range("E8").Value = '  +0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.09363'
$ws.Range("E9").Value = '  -1.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.139'
$ws.Range("E10").Value = '  +3.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '41.75'
$ws.Range("E11").Value = '  +3.18%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.442'
$ws.Range("E12").Value = '  -0.12%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.93'
$ws.Range("E13").Value = '  +1.73%  '
$ws.Range("D14").Value = '1.875.70'
$ws.Range("E14").Value = '  +3.42%  '
$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.380'
$ws.Range("E15").Value = '  +1.27%  '
$ws.Range("B16").Value = 'BinanceUSD'
$ws.Range("C16").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.003'
$ws.Range("E16").Value = '  +0.21%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001120'
$ws.Range("E17").Value = '  -0.81%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '92.34'
$ws.Range("E18").Value = '  -1.22%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06591'
$ws.Range("E19").Value = '  -0.38%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.74'
$ws.Range("E20").Value = '  +3.37%  '
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("E22").Value = '  +3.35%  '
$ws.Range("D23").Value = '28.277.92'
$ws.Range("E23").Value = '  +1.27%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.31'
$ws.Range("E24").Value = '  +1.11%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.273'
$ws.Range("E25").Value = '  +0.99%  '
$ws.Range("E26").Value = '  +5.82%  '
$ws.Range("D27").Value = '2.090.55'
$ws.Range("E27").Value = '  +3.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '21.12'
$ws.Range("E28").Value = '  +1.24%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '158.99'
$ws.Range("E29").Value = '  +1.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '127.11'
$ws.Range("E30").Value = '  -1.04%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.1058'
$ws.Range("E31").Value = '  -0.94%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.066'
$ws.Range("E32").Value = '  +1.10%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.613'
$ws.Range("E33").Value = '  +0.35%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.625'
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.06713'
$ws.Range("E35").Value = '  -1.85%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '9.471'
$ws.Range("E36").Value = '  +4.26%  '
$ws.Range("E37").Value = '  +4.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2181'
$ws.Range("E38").Value = '  +1.04%  '
$ws.Range("B39").Value = 'InternetComputer(DFINITY)'
$ws.Range("C39").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.035'
$ws.Range("E39").Value = '  +1.37%  '
$ws.Range("B40").Value = 'Aptos'
$ws.Range("C40").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '11.45'
$ws.Range("E40").Value = '  +0.70%  '
$ws.Range("B41").Value = 'TheSandbox'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6342'
$ws.Range("E41").Value = '  +1.86%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.188'
$ws.Range("E42").Value = '  +3.44%  '
$ws.Range("E43").Value = '  +0.08%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '13.48'
$ws.Range("E44").Value = '  +2.64%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5964'
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("E46").Value = '  -1.62%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.652'
$ws.Range("E47").Value = '  -1.39%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.998'
$ws.Range("E48").Value = '  +1.97%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '123.37'
$ws.Range("E49").Value = '  -0.83%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.195'
$ws.Range("E50").Value = '  +1.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06857'
$ws.Range("E51").Value = '  +0.92%  '
